$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute('2026-01-26 Monday', $true, $false, $false, $false, $false, $true, 1, $false, '2026-01-27 Tuesday', 2) | Out-Null

# New values for each of the 100 table cells, in row-major order
$newValues = @(
    '23-10=',
    '34+2=',
    '38+34=',
    '65+33=',
    '55+10=',
    '61-53=',
    '5+34=',
    '42+23=',
    '35+20=',
    '94-59=',
    '52-48=',
    '14-1=',
    '72-48=',
    '33+52=',
    '60+37=',
    '34+3=',
    '82-66=',
    '54+34=',
    '9+9=',
    '20+17=',
    '19+9=',
    '34-8=',
    '71+21=',
    '50-16=',
    '17+36=',
    '52-20=',
    '96-8=',
    '62+31=',
    '12+66=',
    '44-40=',
    '12+74=',
    '70+11=',
    '85+5=',
    '35+20=',
    '41+5=',
    '39+50=',
    '25+17=',
    '52-22=',
    '86-28=',
    '50+5=',
    '94+0=',
    '34+32=',
    '20-15=',
    '48-33=',
    '86-1=',
    '30+20=',
    '44+41=',
    '73-41=',
    '68+29=',
    '7+0=',
    '97-53=',
    '80-60=',
    '57-27=',
    '87+10=',
    '59-12=',
    '22-17=',
    '11+1=',
    '11+88=',
    '51-11=',
    '54+45=',
    '56+14=',
    '36-33=',
    '13+4=',
    '13+30=',
    '91-16=',
    '70-10=',
    '19+30=',
    '37+58=',
    '87+8=',
    '34+14=',
    '31+22=',
    '84-34=',
    '5+47=',
    '35+47=',
    '82-69=',
    '32+23=',
    '31+57=',
    '92-3=',
    '37-29=',
    '25+34=',
    '17+34=',
    '22+15=',
    '20+64=',
    '3-0=',
    '7+44=',
    '43+34=',
    '90-0=',
    '28+19=',
    '31-11=',
    '92-11=',
    '48+21=',
    '33-32=',
    '18+68=',
    '49-1=',
    '87-36=',
    '79-70=',
    '43-18=',
    '64+12=',
    '57-23=',
    '71-25='
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output "Updated $idx cells"
